$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: populate the first data row of the report ---

# A7:B7 - row number (was 0, now 1)
$ws.Range("A7").Value = 1

# C7:G7 and N7:O7 share the same direct-formatting style (numFmtId 0 -> 49 / Text)
# in the target workbook, so switch both ranges to Text format before writing the
# text values (this lets the two ranges keep sharing one style, same as before).
$ws.Range("N7:O7").NumberFormat = "@"
$ws.Range("C7:G7").NumberFormat = "@"
$ws.Range("C7").Value = "جونتي عمال"
$ws.Range("N7").Value = "40.00"

# H7:K7 - its own style also flips numFmtId 0 -> 49 / Text
$ws.Range("H7:K7").NumberFormat = "@"
$ws.Range("H7").Value = "5:0"

# L7:M7 keeps its original numeric style (numFmtId 165) but the stored cell
# becomes a text value - toggle to Text, write the value, then restore the
# original custom format so the style index itself is unchanged.
$ws.Range("L7:M7").NumberFormat = "@"
$ws.Range("L7").Value = "0"
$ws.Range("L7:M7").NumberFormat = "#,##0.##;""[""#,##0.##""]"";0"

# P7 keeps its original numeric style (numFmtId 2, i.e. "0.00") but also becomes
# a text value the same way.
$ws.Range("P7").NumberFormat = "@"
$ws.Range("P7").Value = "40.0000"
$ws.Range("P7").NumberFormat = "0.00"

# Q7 - own style flips numFmtId 0 -> 49 / Text
$ws.Range("Q7").NumberFormat = "@"
$ws.Range("Q7").Value = "1:0"

# --- Row 8: totals row ---
$ws.Range("P8").Value = 40

# --- Row 9: footer / generated timestamp ---
$ws.Range("A9").Value = "Sunday, 21 September, 2025 9:25 AM"

Write-Output "edit applied"
